# Update revenue projections (Clinton column) to reflect Clinton CTC update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = 18.406863006094078
$ws.Range("C4").Value = 18.160015296595475
$ws.Range("C5").Value = 18.517215569811132
$ws.Range("C6").Value = 18.501408786788229
$ws.Range("C7").Value = 18.51288048582191
$ws.Range("C8").Value = 18.626240949857507
$ws.Range("C9").Value = 18.575050498562437
$ws.Range("C10").Value = 18.625631473037672
$ws.Range("C11").Value = 18.645986046475375
$ws.Range("C12").Value = 18.707127769652164
$ws.Range("C13").Value = 18.801497904249576
$ws.Range("C14").Value = 18.868966720585707
$ws.Range("C15").Value = 18.851037852570972
$ws.Range("C16").Value = 18.970076554976252
$ws.Range("C17").Value = 19.054455820249423
$ws.Range("C18").Value = 19.051898093245796
$ws.Range("C19").Value = 19.20556643427738
$ws.Range("C20").Value = 19.169319891986099
$ws.Range("C21").Value = 19.325287829431627
$ws.Range("C22").Value = 19.40941223991031
$ws.Range("C23").Value = 19.404358430568255
$ws.Range("C24").Value = 19.539252365850686
$ws.Range("C25").Value = 19.628883722401753
$ws.Range("C26").Value = 19.631312428188032
$ws.Range("C27").Value = 19.761607379743378
$ws.Range("C28").Value = 19.739861281698794

# Update active selection to B1:E1 (header row)
$ws.Range("B1:E1").Select() | Out-Null
